$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New project rows to append below the existing data (A5:B26).
# Column A = project code, Column B = budget ("N/A" or a numeric-looking
# amount such as "15000.00" that must be stored as *text*, matching the
# rest of the budget column).
$data = @(
    @("AL071236", "N/A"),
    @("AL112347", "N/A"),
    @("AL213557", "N/A"),
    @("AL221556", "N/A"),
    @("AL232310", "N/A"),
    @("AL258112", "N/A"),
    @("AL258113", "N/A"),
    @("AR044321", "N/A"),
    @("AR056691", "15000.00"),
    @("AR062233", "N/A"),
    @("AR062331", "N/A"),
    @("UL046321", "17000.00"),
    @("UL054239", "17000.00"),
    @("UL223314", "N/A"),
    @("UL258001", "N/A"),
    @("UL258122", "N/A"),
    @("UL258129", "N/A"),
    @("UR047451", "14166.67"),
    @("UR054912", "N/A"),
    @("UR216878", "N/A"),
    @("UR237511", "N/A"),
    @("UR332441", "N/A")
)

$row = 5
foreach ($pair in $data) {
    $code = $pair[0]
    $budget = $pair[1]

    $ws.Cells.Item($row, 1).Value = $code

    $budgetCell = $ws.Cells.Item($row, 2)
    if ($budget -match '^[0-9]+(\.[0-9]+)?$') {
        # Numeric-looking text (e.g. "15000.00") must be preserved as a
        # literal text string rather than auto-converted to a number.
        # Build it with TEXT() and then "bake" the formula result into a
        # plain value via copy / paste-values, so the stored cell keeps
        # the General number format (no new style is introduced) while
        # its content is stored as text.
        $budgetCell.Formula = '=TEXT(' + $budget + ',"0.00")'
        $budgetCell.Copy()
        $budgetCell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
    } else {
        $budgetCell.Value = $budget
    }

    $row = $row + 1
}

$excel.CutCopyMode = $false
